$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = "new"
$ws.Cells.Item(2, 3).Value = 653
$ws.Cells.Item(3, 2).Value = "vaccination"
$ws.Cells.Item(3, 3).Value = 638
$ws.Cells.Item(4, 2).Value = "incidence"
$ws.Cells.Item(4, 3).Value = 605
$ws.Cells.Item(5, 2).Value = "new infections"
$ws.Cells.Item(5, 3).Value = 605
$ws.Cells.Item(6, 2).Value = "infected"
$ws.Cells.Item(6, 3).Value = 410
$ws.Cells.Item(7, 2).Value = "people"
$ws.Cells.Item(7, 3).Value = 401
$ws.Cells.Item(8, 2).Value = "persons"
$ws.Cells.Item(8, 3).Value = 347
$ws.Cells.Item(9, 2).Value = "deaths"
$ws.Cells.Item(9, 3).Value = 320
$ws.Cells.Item(10, 2).Value = "pandemic"
$ws.Cells.Item(10, 3).Value = 319
$ws.Cells.Item(11, 2).Value = "compulsory vaccination"
$ws.Cells.Item(11, 3).Value = 285
$ws.Cells.Item(12, 2).Value = "vaccination"
$ws.Cells.Item(12, 3).Value = 455
$ws.Cells.Item(13, 2).Value = "incidence"
$ws.Cells.Item(13, 3).Value = 454
$ws.Cells.Item(14, 2).Value = "new"
$ws.Cells.Item(14, 3).Value = 412
$ws.Cells.Item(15, 2).Value = "new infections"
$ws.Cells.Item(15, 3).Value = 382
$ws.Cells.Item(16, 2).Value = "people"
$ws.Cells.Item(16, 3).Value = 306
$ws.Cells.Item(17, 2).Value = "infected"
$ws.Cells.Item(17, 3).Value = 304
$ws.Cells.Item(18, 2).Value = "deaths"
$ws.Cells.Item(18, 3).Value = 268
$ws.Cells.Item(19, 2).Value = "persons"
$ws.Cells.Item(19, 3).Value = 256
$ws.Cells.Item(20, 2).Value = "pandemic"
$ws.Cells.Item(20, 3).Value = 215
$ws.Cells.Item(21, 2).Value = "gives"
$ws.Cells.Item(21, 3).Value = 208
$ws.Cells.Item(22, 2).Value = "incidence"
$ws.Cells.Item(22, 3).Value = 449
$ws.Cells.Item(23, 2).Value = "new infections"
$ws.Cells.Item(23, 3).Value = 413
$ws.Cells.Item(24, 2).Value = "vaccination"
$ws.Cells.Item(24, 3).Value = 371
$ws.Cells.Item(25, 2).Value = "new"
$ws.Cells.Item(25, 3).Value = 335
$ws.Cells.Item(26, 2).Value = "people"
$ws.Cells.Item(26, 3).Value = 309
$ws.Cells.Item(27, 2).Value = "persons"
$ws.Cells.Item(27, 3).Value = 285
$ws.Cells.Item(28, 2).Value = "pandemic"
$ws.Cells.Item(28, 3).Value = 279
$ws.Cells.Item(29, 2).Value = "infected"
$ws.Cells.Item(29, 3).Value = 232
$ws.Cells.Item(30, 2).Value = "deaths"
$ws.Cells.Item(30, 3).Value = 229
$ws.Cells.Item(31, 2).Value = "gives"
$ws.Cells.Item(31, 3).Value = 228
$ws.Cells.Item(32, 2).Value = "incidence"
$ws.Cells.Item(32, 3).Value = 284
$ws.Cells.Item(33, 2).Value = "vaccination"
$ws.Cells.Item(33, 3).Value = 254
$ws.Cells.Item(34, 2).Value = "new infections"
$ws.Cells.Item(34, 3).Value = 237
$ws.Cells.Item(35, 2).Value = "people"
$ws.Cells.Item(35, 3).Value = 195
$ws.Cells.Item(36, 2).Value = "new"
$ws.Cells.Item(36, 3).Value = 184
$ws.Cells.Item(37, 2).Value = "data"
$ws.Cells.Item(37, 3).Value = 166
$ws.Cells.Item(38, 2).Value = "persons"
$ws.Cells.Item(38, 3).Value = 156
$ws.Cells.Item(39, 2).Value = "deaths"
$ws.Cells.Item(39, 3).Value = 155
$ws.Cells.Item(40, 2).Value = "infections"
$ws.Cells.Item(40, 3).Value = 151
$ws.Cells.Item(41, 2).Value = "pandemic"
$ws.Cells.Item(41, 3).Value = 139
$ws.Cells.Item(42, 2).Value = "new infections"
$ws.Cells.Item(42, 3).Value = 214
$ws.Cells.Item(43, 2).Value = "vaccination"
$ws.Cells.Item(43, 3).Value = 201
$ws.Cells.Item(44, 2).Value = "incidence"
$ws.Cells.Item(44, 3).Value = 180
$ws.Cells.Item(45, 2).Value = "new"
$ws.Cells.Item(45, 3).Value = 176
$ws.Cells.Item(46, 2).Value = "people"
$ws.Cells.Item(46, 3).Value = 176
$ws.Cells.Item(47, 2).Value = "monkeypox"
$ws.Cells.Item(47, 3).Value = 158
$ws.Cells.Item(48, 2).Value = "deaths"
$ws.Cells.Item(48, 3).Value = 133
$ws.Cells.Item(49, 2).Value = "infections"
$ws.Cells.Item(49, 3).Value = 119
$ws.Cells.Item(50, 2).Value = "pandemic"
$ws.Cells.Item(50, 3).Value = 115
$ws.Cells.Item(51, 2).Value = "data"
$ws.Cells.Item(51, 3).Value = 110
$ws.Cells.Item(52, 2).Value = "vaccination"
$ws.Cells.Item(52, 3).Value = 212
$ws.Cells.Item(53, 2).Value = "new infections"
$ws.Cells.Item(53, 3).Value = 180
$ws.Cells.Item(54, 2).Value = "new"
$ws.Cells.Item(54, 3).Value = 156
$ws.Cells.Item(55, 2).Value = "people"
$ws.Cells.Item(55, 3).Value = 150
$ws.Cells.Item(56, 2).Value = "incidence"
$ws.Cells.Item(56, 3).Value = 145
$ws.Cells.Item(57, 2).Value = "gives"
$ws.Cells.Item(57, 3).Value = 135
$ws.Cells.Item(58, 2).Value = "autumn"
$ws.Cells.Item(58, 3).Value = 114
$ws.Cells.Item(59, 2).Value = "children"
$ws.Cells.Item(59, 3).Value = 113
$ws.Cells.Item(60, 2).Value = "pandemic"
$ws.Cells.Item(60, 3).Value = 100
$ws.Cells.Item(61, 2).Value = "pay"
$ws.Cells.Item(61, 3).Value = 93
$ws.Cells.Item(62, 2).Value = "vaccination"
$ws.Cells.Item(62, 3).Value = 228
$ws.Cells.Item(63, 2).Value = "new infections"
$ws.Cells.Item(63, 3).Value = 163
$ws.Cells.Item(64, 2).Value = "new"
$ws.Cells.Item(64, 3).Value = 141
$ws.Cells.Item(65, 2).Value = "people"
$ws.Cells.Item(65, 3).Value = 140
$ws.Cells.Item(66, 2).Value = "incidence"
$ws.Cells.Item(66, 3).Value = 119
$ws.Cells.Item(67, 2).Value = "gives"
$ws.Cells.Item(67, 3).Value = 112
$ws.Cells.Item(68, 2).Value = "infection"
$ws.Cells.Item(68, 3).Value = 94
$ws.Cells.Item(69, 2).Value = "pandemic"
$ws.Cells.Item(69, 3).Value = 90
$ws.Cells.Item(70, 2).Value = "7-day #incidence"
$ws.Cells.Item(70, 3).Value = 87
$ws.Cells.Item(71, 2).Value = "infected"
$ws.Cells.Item(71, 3).Value = 82
$ws.Cells.Item(72, 2).Value = "vaccination"
$ws.Cells.Item(72, 3).Value = 176
$ws.Cells.Item(73, 2).Value = "new infections"
$ws.Cells.Item(73, 3).Value = 140
$ws.Cells.Item(74, 2).Value = "infection"
$ws.Cells.Item(74, 3).Value = 115
$ws.Cells.Item(75, 2).Value = "new"
$ws.Cells.Item(75, 3).Value = 114
$ws.Cells.Item(76, 2).Value = "incidence"
$ws.Cells.Item(76, 3).Value = 106
$ws.Cells.Item(77, 2).Value = "gives"
$ws.Cells.Item(77, 3).Value = 94
$ws.Cells.Item(78, 2).Value = "7-day #incidence"
$ws.Cells.Item(78, 3).Value = 94
$ws.Cells.Item(79, 2).Value = "people"
$ws.Cells.Item(79, 3).Value = 92
$ws.Cells.Item(80, 2).Value = "deaths"
$ws.Cells.Item(80, 3).Value = 73
$ws.Cells.Item(81, 2).Value = "study"
$ws.Cells.Item(81, 3).Value = 71
$ws.Cells.Item(82, 2).Value = "vaccination"
$ws.Cells.Item(82, 3).Value = 153
$ws.Cells.Item(83, 2).Value = "new"
$ws.Cells.Item(83, 3).Value = 107
$ws.Cells.Item(84, 2).Value = "new infections"
$ws.Cells.Item(84, 3).Value = 95
$ws.Cells.Item(85, 2).Value = "incidence"
$ws.Cells.Item(85, 3).Value = 87
$ws.Cells.Item(86, 2).Value = "pandemic"
$ws.Cells.Item(86, 3).Value = 79
$ws.Cells.Item(87, 2).Value = "gives"
$ws.Cells.Item(87, 3).Value = 77
$ws.Cells.Item(88, 2).Value = "7-day #incidence"
$ws.Cells.Item(88, 3).Value = 76
$ws.Cells.Item(89, 2).Value = "people"
$ws.Cells.Item(89, 3).Value = 67
$ws.Cells.Item(90, 2).Value = "Friday"
$ws.Cells.Item(90, 3).Value = 63
$ws.Cells.Item(91, 2).Value = "new"
$ws.Cells.Item(91, 3).Value = 62
$ws.Cells.Item(92, 2).Value = "vaccination"
$ws.Cells.Item(92, 3).Value = 182
$ws.Cells.Item(93, 2).Value = "new"
$ws.Cells.Item(93, 3).Value = 132
$ws.Cells.Item(94, 2).Value = "new infections"
$ws.Cells.Item(94, 3).Value = 129
$ws.Cells.Item(95, 2).Value = "people"
$ws.Cells.Item(95, 3).Value = 107
$ws.Cells.Item(96, 2).Value = "incidence"
$ws.Cells.Item(96, 3).Value = 97
$ws.Cells.Item(97, 2).Value = "7-day #incidence"
$ws.Cells.Item(97, 3).Value = 79
$ws.Cells.Item(98, 2).Value = "pandemic"
$ws.Cells.Item(98, 3).Value = 73
$ws.Cells.Item(99, 2).Value = "gives"
$ws.Cells.Item(99, 3).Value = 72
$ws.Cells.Item(100, 2).Value = "population"
$ws.Cells.Item(100, 3).Value = 71
$ws.Cells.Item(101, 2).Value = "antibody"
$ws.Cells.Item(101, 3).Value = 66
$ws.Cells.Item(102, 2).Value = "vaccination"
$ws.Cells.Item(102, 3).Value = 74
$ws.Cells.Item(103, 2).Value = "children"
$ws.Cells.Item(103, 3).Value = 51
$ws.Cells.Item(104, 2).Value = "people"
$ws.Cells.Item(104, 3).Value = 49
$ws.Cells.Item(105, 2).Value = "new infections"
$ws.Cells.Item(105, 3).Value = 47
$ws.Cells.Item(106, 2).Value = "incidence"
$ws.Cells.Item(106, 3).Value = 35
$ws.Cells.Item(107, 2).Value = "Thursday"
$ws.Cells.Item(107, 3).Value = 34
$ws.Cells.Item(108, 2).Value = "study"
$ws.Cells.Item(108, 3).Value = 33
$ws.Cells.Item(109, 2).Value = "new"
$ws.Cells.Item(109, 3).Value = 32
$ws.Cells.Item(110, 2).Value = "infection"
$ws.Cells.Item(110, 3).Value = 29
$ws.Cells.Item(111, 2).Value = "7-day #incidence"
$ws.Cells.Item(111, 3).Value = 25
